$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update "Mio Algoritmo" (C column) results with the new simulated-annealing
# run's values. The "%" (D column) and "Media errori" (D12) formulas will
# recalculate automatically.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 6110
$ws.Range("C3").Value = 15781
$ws.Range("C4").Value = 538
$ws.Range("C5").Value = 22869
$ws.Range("C6").Value = 21282
$ws.Range("C7").Value = 42029
$ws.Range("C8").Value = 51105
$ws.Range("C9").Value = 107330
$ws.Range("C10").Value = 9222
$ws.Range("C11").Value = 231383

# ---------------------------------------------------------------------------
# New "SEED" column (E): header + per-run seed values.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "SEED"
$ws.Range("E1").Interior.Color = 65535
$ws.Range("E1").Borders.Item(7).LineStyle = 1
$ws.Range("E1").Borders.Item(10).LineStyle = 1

$ws.Range("D2").Copy()
$ws.Range("E2:E11").PasteSpecial(-4122)
$ws.Range("E2:E11").NumberFormat = "0"

$ws.Range("E2").Value = 1553592951193
$ws.Range("E3").Value = 1553593210405
$ws.Range("E4").Value = 1553593570081
$ws.Range("E5").Value = 1553793401381
$ws.Range("E6").Value = 1553593959586
$ws.Range("E7").Value = 1553681533325
$ws.Range("E8").Value = 1553707441798
$ws.Range("E9").Value = 1553595355487
$ws.Range("E10").Value = 1553723496331
$ws.Range("E11").Value = 1553610445037

$ws.Columns.Item(5).AutoFit()

# ---------------------------------------------------------------------------
# Status indicator cells (F column) - manual "traffic light" colouring used
# while tuning the simulated-annealing run.
# ---------------------------------------------------------------------------
$ws.Range("F2").Interior.Color = 5287936
$ws.Range("F2").Font.Color = 5287936

$ws.Range("F3").Interior.Color = 49407
$ws.Range("F4").Interior.Color = 5287936
$ws.Range("F5").Interior.Color = 255
$ws.Range("F6").Interior.Color = 5287936
$ws.Range("F7").Interior.Color = 49407

# ---------------------------------------------------------------------------
# Selection left where the author was last working.
# ---------------------------------------------------------------------------
$ws.Range("C10").Select()
